# Update "想去人数" (want-to-go count) figures to the newly generated values.
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 75
$ws1.Range("F3").Value = 11793
$ws1.Range("F7").Value = 11734
$ws1.Range("F11").Value = 37
$ws1.Range("F13").Value = 5808

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 573

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 573
$ws4.Range("F3").Value = 75
$ws4.Range("F5").Value = 11793
$ws4.Range("F9").Value = 11734
$ws4.Range("F13").Value = 37
$ws4.Range("F16").Value = 5808
